# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the other header cells (copied from G1), and fill
# in the values for rows 2-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1, style index 1:
# bold font, border, centered alignment) onto the new header cell H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column data.
$saveValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 1
    8 = 1
    9 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}

$excel.CutCopyMode = 0
